# Fixed Bento 80 Test scripts
#
# The CasesTab / SamplesTab / FilesTab Cypher queries on the "startup" sheet
# (columns B2, B3, B4) each gain a trailing "order By ... LIMIT 100" clause.
# Single-quoted here-strings are used below so the literal backticks
# (Cypher `Column Name` quoting), dollar signs and tabs in the query text are
# not touched by PowerShell escape/expansion rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: CasesTab query ---
$casesQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss:study_subject)
	WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
RETURN DISTINCT 
	ss.study_subject_id AS `Case ID`,
	p.program_acronym AS `Program Code`,
	p.program_id AS `Program ID`,
	s.study_acronym AS `Arm`,
	ss.disease_subtype AS `Diagnosis`,
	sf.grouped_recurrence_score AS `Recurrence Score`,
	d.tumor_size_group AS `Tumor Size (cm)`,
	d.er_status AS `ER Status`,
	d.pr_status AS `PR Status`,
	demo.age_at_index AS `Age (years)`,
	demo.survival_time AS `Survival (days)`
 order By ss.study_subject_id ASC LIMIT 100 
'@

# Trim only the single trailing newline the here-string terminator adds; the
# trailing space after "LIMIT 100" must survive.
$casesQuery = $casesQuery -replace "`r?`n$", ""
$ws.Range("B2").Value2 = $casesQuery

# --- B3: SamplesTab query ---
$samplesQuery = @'
MATCH (ss:study_subject)
	WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[*..2]-(parent)<--(f:file)
OPTIONAL MATCH (f)-[:file_of_laboratory_procedure]->(lp)
RETURN DISTINCT 
	samp.sample_id AS `Sample ID`,
	ss.study_subject_id AS `Case ID`,
	p.program_acronym AS `Program Code`,
	s.study_acronym AS `Arm`,
	ss.disease_subtype AS `Diagnosis`,
	samp.tissue_type AS `Tissue Type`,
	samp.composition AS `Tissue Composition`,
	samp.sample_anatomic_site AS `Sample Anatomic Site`,
	samp.method_of_sample_procurement AS `Sample Procurement Method`,
	lp.test_name AS `platform`
 order By samp.sample_id ASC LIMIT 100
'@
$samplesQuery = $samplesQuery -replace "`r?`n$", ""
$ws.Range("B3").Value2 = $samplesQuery

# --- B4: FilesTab query ---
$filesQuery = @'
MATCH (ss:study_subject)
	WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (ss)<-[*..2]-(parent)<--(f:file)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
RETURN DISTINCT 
	f.file_name AS `File Name`,
	head(labels(parent)) AS `Association`,
	f.file_description AS `Description`,
	f.file_format AS `File Format`,
	f.file_size AS `Size`,
	p.program_acronym AS `Program Code`,
	s.study_acronym AS `Arm`,
	ss.study_subject_id AS `Case ID`,
	samp.sample_id AS `Sample ID`
 order By f.file_name ASC LIMIT 100
'@
$filesQuery = $filesQuery -replace "`r?`n$", ""
$ws.Range("B4").Value2 = $filesQuery

# --- Row heights grow slightly to fit the extra wrapped line of text ---
$ws.Rows.Item(2).RowHeight = 345.6
$ws.Rows.Item(3).RowHeight = 259.2
$ws.Rows.Item(4).RowHeight = 244.8

# --- Selection moved from D4:E4 to C4 ---
$ws.Range("C4").Select()
